$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row to add (row 13), matching the style of existing data rows
$newRow = 13

# Column A: date value 45917 (2025-09-17), formatted like other date cells (style index 2 / "YYYY-MM-DD HH:MM:SS")
$ws.Cells.Item($newRow, 1).Value = 45917
$ws.Cells.Item($newRow, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Columns B-E: text values (comma decimal separator), stored as text/inline strings
$ws.Cells.Item($newRow, 2).Value = "20,9899"
$ws.Cells.Item($newRow, 3).Value = "10,7342"
$ws.Cells.Item($newRow, 4).Value = "14,8091"
$ws.Cells.Item($newRow, 5).Value = "14,8091"
